# Update gh-pages to output generated at 456a3b4
# Applies updated "want-to-go" counts (column F) across the four sheets,
# plus one refreshed cover image link (I31) on sheet "展览".

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value  = 318
$ws1.Range("F5").Value  = 5733
$ws1.Range("F7").Value  = 9740
$ws1.Range("F8").Value  = 69
$ws1.Range("F10").Value = 3887
$ws1.Range("F14").Value = 209
$ws1.Range("F18").Value = 110
$ws1.Range("F20").Value = 623
$ws1.Range("F21").Value = 3916
$ws1.Range("F22").Value = 138
$ws1.Range("F24").Value = 5370
$ws1.Range("F26").Value = 2124
$ws1.Range("F27").Value = 135
$ws1.Range("F28").Value = 360
$ws1.Range("F29").Value = 7992
$ws1.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202407/vAqeEX2p1722411672295.png"
$ws1.Range("F32").Value = 2207
$ws1.Range("F33").Value = 2213
$ws1.Range("F35").Value = 1313
$ws1.Range("F38").Value = 273
$ws1.Range("F39").Value = 251
$ws1.Range("F40").Value = 20
$ws1.Range("F41").Value = 1187
$ws1.Range("F42").Value = 1180
$ws1.Range("F44").Value = 1345
$ws1.Range("F45").Value = 2113
$ws1.Range("F46").Value = 135
$ws1.Range("F47").Value = 229

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F7").Value  = 2
$ws2.Range("F9").Value  = 935
$ws2.Range("F15").Value = 13
$ws2.Range("F20").Value = 17

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 583
$ws3.Range("F3").Value = 761
$ws3.Range("F4").Value = 69

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 761
$ws4.Range("F4").Value  = 69
$ws4.Range("F5").Value  = 318
$ws4.Range("F6").Value  = 5733
$ws4.Range("F8").Value  = 3887
$ws4.Range("F15").Value = 110
$ws4.Range("F18").Value = 623
$ws4.Range("F19").Value = 3916
$ws4.Range("F21").Value = 138
$ws4.Range("F22").Value = 2
$ws4.Range("F23").Value = 5370
$ws4.Range("F25").Value = 2124
$ws4.Range("F26").Value = 135
$ws4.Range("F27").Value = 360
$ws4.Range("F28").Value = 7992
$ws4.Range("F30").Value = 2207
$ws4.Range("F31").Value = 2213
$ws4.Range("F33").Value = 1313
$ws4.Range("F35").Value = 273
$ws4.Range("F36").Value = 251
$ws4.Range("F37").Value = 20
$ws4.Range("F38").Value = 1187
$ws4.Range("F39").Value = 1180
$ws4.Range("F42").Value = 1345
$ws4.Range("F43").Value = 13
$ws4.Range("F44").Value = 2113
$ws4.Range("F45").Value = 135
$ws4.Range("F46").Value = 229
$ws4.Range("F48").Value = 17
